$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @(1,0,0,1,1,1,1),
    @(0,1,1,1,0,0,1),
    @(0,1,1,0,1,0,0),
    @(1,1,0,1,1,0,1),
    @(1,0,1,1,1,0,1),
    @(1,0,0,0,0,1,0),
    @(1,1,0,1,1,0,1)
)

for ($r = 1; $r -le 7; $r++) {
    for ($c = 1; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$r - 1][$c - 1]
    }
}

$ws.Range("G7").Select()
